# The diary entry that used to read "02.05.2017 - 03.05.2017" is split into
# "02.05.2017 - 05" + ".05.2017" (two separate runs), the trailing _GoBack
# bookmark is removed from that paragraph, and a batch of new paragraphs is
# appended right after "Das Bestellarray funktioniert jetzt und wurde
# perfektioniert.", with the _GoBack bookmark relocated to the very end of
# the new content.

$d = $word.ActiveDocument

# Locate the two relevant paragraphs by their text instead of a hard-coded
# index, so the script keeps working even if the paragraph numbering shifts.
$dateParaIndex = -1
$bestellarrayParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -like "02.05.2017*03.05.2017*") {
        $dateParaIndex = $i
    }
    if ($text -like "Das Bestellarray funktioniert jetzt*") {
        $bestellarrayParaIndex = $i
    }
}

# Step 1: the _GoBack bookmark currently sits at the end of the
# "02.05.2017 - 03.05.2017" paragraph. It needs to move to the end of the
# new content, so remove it from its current spot first.
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

# Step 2: change " - 03.05.2017" to " - 05" (still a single run) within that
# paragraph, then insert a new run ".05.2017" right after it, producing
# three runs in total for that paragraph.
$dateParagraph = $d.Paragraphs.Item($dateParaIndex)
$dateRange = $d.Range($dateParagraph.Range.Start, $dateParagraph.Range.End)
$dateRange.Find.Execute(" – 03.05.2017", $true, $false, $false, $false, `
    $false, $true, 1, $false, " – 05", 2) | Out-Null

$dateParagraph = $d.Paragraphs.Item($dateParaIndex)
$dateTailPos = $dateParagraph.Range.End - 1
$d.Range($dateTailPos, $dateTailPos).InsertAfter(".05.2017")

# Step 3: append the new paragraphs right after "Das Bestellarray ..." using
# a WordML fragment, so the freshly created empty paragraph stays truly
# empty (<w:p/>) instead of gaining a stray empty run, and the relocated
# bookmark keeps its original id/name.
$bestellarrayParagraph = $d.Paragraphs.Item($bestellarrayParaIndex)
$insertPos = $bestellarrayParagraph.Range.End - 1
$insertionPoint = $d.Range($insertPos, $insertPos)

$newContent = '<?xml version="1.0"?>' + `
  '<w:wordDocument xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + `
      '<w:p></w:p>' + `
      '<w:p><w:r><w:t>05.05.2017</w:t></w:r></w:p>' + `
      '<w:p><w:r><w:t>Liste angefangen</w:t></w:r></w:p>' + `
      '<w:p><w:r><w:t>08.05.2017</w:t></w:r></w:p>' + `
      '<w:p><w:r><w:t>Liste fertig, wird aber durch eine Tabelle ersetzt</w:t></w:r>' + `
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
    '</w:body>' + `
  '</w:wordDocument>'

$insertionPoint.InsertXML($newContent)
